$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2-9 (rank table). A new entry is being
# inserted as the new 4th place ("testando" / 1919), which pushes the
# former rows 5-9 down by one (to rows 6-10). Two brand new entries are
# then appended at the bottom as rows 11 and 12.

# Make sure any numeric-looking text we write into column D stays TEXT
# (matching the existing cells in that column) instead of being
# auto-converted to a number.
$ws.Range("D5:D12").NumberFormat = "@"

# --- Row 5: new 4th-place entry ---
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "4º"
$ws.Cells.Item(5, 3).Value = "testando"
$ws.Cells.Item(5, 4).Value = "1919"

# --- Rows 6-10: former rows 5-9 shifted down one position ---
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "5º"
$ws.Cells.Item(6, 3).Value = "diego"
$ws.Cells.Item(6, 4).Value = "1891"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "6º"
$ws.Cells.Item(7, 3).Value = "diego"
$ws.Cells.Item(7, 4).Value = "1826"

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "7º"
$ws.Cells.Item(8, 3).Value = "diego"
$ws.Cells.Item(8, 4).Value = "1791"

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "8º"
$ws.Cells.Item(9, 3).Value = "diego"
$ws.Cells.Item(9, 4).Value = "1424"

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "9º"
$ws.Cells.Item(10, 3).Value = "diego"
$ws.Cells.Item(10, 4).Value = "1423"

# New row 10/11/12 column-A cells need the same bold/bordered style that
# the rest of column A already uses (rows 2-9 have it); copy it down.
$ws.Cells.Item(9, 1).Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 11: brand new entry ---
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "10º"
$ws.Cells.Item(11, 3).Value = "testando2"
$ws.Cells.Item(11, 4).Value = "0"

# --- Row 12: brand new entry ---
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "11º"
$ws.Cells.Item(12, 3).Value = "teatando2"
$ws.Cells.Item(12, 4).Value = "-500"
